$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (54 and 55) to the feed log sheet
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 1
$ws.Range("C54").Value = "2024-06-16 01:01:55"
$ws.Range("D54").Value = 200
$ws.Range("E54").Value = 8

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 2
$ws.Range("C55").Value = "2024-06-16 01:01:55"
$ws.Range("D55").Value = 200
$ws.Range("E55").Value = 0
